# Re-applies the latest coinranking.com snapshot (prices + 1h volume%) onto
# the cryptos worksheet, including the row-43/44 and row-46/47 rank swaps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Excel auto-converts plain numeric-looking strings (e.g. "1.00") to
    # Double when assigned via .Value, which would silently drop the text
    # formatting the source data relies on (trailing zeros, thousand dots,
    # etc). A leading apostrophe forces text entry like a user typing it in;
    # resetting the style back to Normal afterwards avoids leaving a stray
    # quote-prefixed / text-formatted cell style behind.
    $cell = $ws.Range($range)
    $isNumericLooking = $text -match '^-?\d+(\.\d+)?$'
    if ($isNumericLooking) {
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue "D2" "80.365.60"
Set-TextValue "E2" "  +4.66%  "
Set-TextValue "D3" "3.184.60"
Set-TextValue "E3" "  +1.27%  "
Set-TextValue "E4" "  +0.18%  "
Set-TextValue "D5" "210.06"
Set-TextValue "E5" "  +4.07%  "
Set-TextValue "D6" "627.37"
Set-TextValue "E6" "  +0.12%  "
Set-TextValue "D7" "0.273"
Set-TextValue "E7" "  +25.67%  "
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  +0.10%  "
Set-TextValue "D9" "0.588"
Set-TextValue "E9" "  +5.04%  "
Set-TextValue "D10" "3.183.84"
Set-TextValue "E10" "  +1.21%  "
Set-TextValue "D11" "0.589"
Set-TextValue "E11" "  +23.01%  "
Set-TextValue "E12" "  +27.26%  "
Set-TextValue "E13" "  +1.28%  "
Set-TextValue "D14" "3.774.78"
Set-TextValue "E14" "  +1.71%  "
Set-TextValue "D15" "5.28"
Set-TextValue "E15" "  -0.10%  "
Set-TextValue "D16" "31.93"
Set-TextValue "E16" "  +7.18%  "
Set-TextValue "D17" "80.462.88"
Set-TextValue "E17" "  +4.97%  "
Set-TextValue "D18" "3.192.45"
Set-TextValue "E18" "  +1.70%  "
Set-TextValue "D19" "14.24"
Set-TextValue "E19" "  +3.58%  "
Set-TextValue "D20" "3.02"
Set-TextValue "E20" "  +9.48%  "
Set-TextValue "D21" "9.17"
Set-TextValue "E21" "  -1.34%  "
Set-TextValue "D22" "437.93"
Set-TextValue "E22" "  +9.19%  "
Set-TextValue "D23" "5.21"
Set-TextValue "E23" "  +13.62%  "
Set-TextValue "E24" "  +6.10%  "
Set-TextValue "D25" "3.351.51"
Set-TextValue "E25" "  +1.70%  "
Set-TextValue "D26" "76.19"
Set-TextValue "E26" "  +3.20%  "
Set-TextValue "D27" "4.70"
Set-TextValue "E27" "  +1.33%  "
Set-TextValue "D28" "10.91"
Set-TextValue "E28" "  +5.07%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  -0.01%  "
Set-TextValue "D30" "0.0000122"
Set-TextValue "E30" "  +8.15%  "
Set-TextValue "D31" "0.997"
Set-TextValue "E31" "  +0.34%  "
Set-TextValue "D32" "8.95"
Set-TextValue "E32" "  +5.23%  "
Set-TextValue "D33" "559.50"
Set-TextValue "E33" "  +6.67%  "
Set-TextValue "E34" "  -0.76%  "
Set-TextValue "D35" "0.151"
Set-TextValue "E35" "  +12.29%  "
Set-TextValue "D36" "2.00"
Set-TextValue "E36" "  +1.92%  "
Set-TextValue "D37" "23.05"
Set-TextValue "E37" "  +5.54%  "
Set-TextValue "E38" "  +18.93%  "
Set-TextValue "D39" "1.00"
Set-TextValue "E39" "  +0.05%  "
Set-TextValue "E40" "  +5.12%  "
Set-TextValue "D41" "20.77"
Set-TextValue "E41" "  +3.49%  "
Set-TextValue "D42" "163.08"
Set-TextValue "E42" "  -0.25%  "
Set-TextValue "B43" "RenderToken"
Set-TextValue "C43" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D43" "5.64"
Set-TextValue "E43" "  +5.13%  "
Set-TextValue "B44" "USDe"
Set-TextValue "C44" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D44" "1.00"
Set-TextValue "E44" "  +0.01%  "
Set-TextValue "D45" "189.17"
Set-TextValue "E45" "  -3.61%  "
Set-TextValue "B46" "dogwifhat"
Set-TextValue "C46" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D46" "2.72"
Set-TextValue "E46" "  +9.26%  "
Set-TextValue "B47" "Stacks"
Set-TextValue "C47" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "1.81"
Set-TextValue "E47" "  +5.12%  "
Set-TextValue "D48" "0.783"
Set-TextValue "E48" "  -3.48%  "
Set-TextValue "E49" "  +0.64%  "
Set-TextValue "D50" "42.83"
Set-TextValue "E50" "  +2.21%  "
Set-TextValue "D51" "4.25"
Set-TextValue "E51" "  +5.62%  "
